# Update localization-status report:
#   - "Ready for handoff" -> "In Translation" on all three sheets
#   - shrink the now-narrower Status columns to match the new text length
#
# $excel / $wb are provided by the host. $wb.ActiveWorkbook is already open.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Column width the report recomputes after the text changes. The stored
# OOXML <col width="..."> is ColumnWidth + 5/6 (Excel's "padding" offset),
# so back that out here to land on the closest value the engine can store.
$targetStoredWidth = 13.4101845877511
$targetColumnWidth = 12.576851254417766

# --- Sheet "Overview": zh-cn / de-de status columns are E and F ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns("E:F").ColumnWidth = $targetColumnWidth

# --- Sheet "zh-cn": Status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns("C").ColumnWidth = $targetColumnWidth

# --- Sheet "de-de": Status column is C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns("C").ColumnWidth = $targetColumnWidth
